$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.29"
$ws.Range("E2").Value = "'6.80%"
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3").Value = "'27.40"
$ws.Range("E3").Value = "'3.58%"
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("D4").Value = "'4.804"
$ws.Range("E4").Value = "'1.85%"
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("D5").Value = "'0.06352"
$ws.Range("E5").Value = "'3.14%"
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("D6").Value = "'6.949"
$ws.Range("E6").Value = "'3.58%"
$ws.Range("D6:E6").Style = "Normal"

$ws.Range("D7").Value = "'3.398"
$ws.Range("E7").Value = "'7.08%"
$ws.Range("D7:E7").Style = "Normal"

$ws.Range("E8").Value = "'3.44%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9565"
$ws.Range("E9").Value = "'4.83%"
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1474"
$ws.Range("E10").Value = "'4.79%"
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.05149"
$ws.Range("E11").Value = "'-3.08%"
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07315"
$ws.Range("E12").Value = "'3.03%"
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03121"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("D13:E13").Style = "Normal"

$ws.Range("D14").Value = "'0.09072"
$ws.Range("E14").Value = "'0.30%"
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001564"
$ws.Range("E15").Value = "'1.57%"
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'0.0006266"
$ws.Range("E16").Value = "'1.38%"
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'0.005797"
$ws.Range("E17").Value = "'-2.78%"
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18").Value = "'3.453"
$ws.Range("E18").Value = "'0.06%"
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("E19").Value = "'4.97%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D22").Value = "'3.871"
$ws.Range("E22").Value = "'-5.25%"
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("E23").Value = "'2.04%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001177"
$ws.Range("E24").Value = "'-0.22%"
$ws.Range("D24:E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004291"
$ws.Range("E25").Value = "'6.00%"
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("E26").Value = "'0.00%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0001692"
$ws.Range("E27").Value = "'3.15%"
$ws.Range("D27:E27").Style = "Normal"

$ws.Range("D40").Value = "'0.04089"
$ws.Range("E40").Value = "'2.49%"
$ws.Range("D40:E40").Style = "Normal"

$ws.Range("D41").Value = "'0.006666"
$ws.Range("E41").Value = "'61.92%"
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("E42").Value = "'4.65%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002202"
$ws.Range("E43").Value = "'2.75%"
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01256"
$ws.Range("E44").Value = "'-5.51%"
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005218"
$ws.Range("E45").Value = "'1.06%"
$ws.Range("D45:E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("D46:E46").Style = "Normal"

$ws.Range("E47").Value = "'821.76%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.02252"
$ws.Range("E48").Value = "'6.16%"
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("D49:E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("D50:E50").Style = "Normal"

